# Update column A (rows 1-100) on the active worksheet with new search-term values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'womens freezer suit',
    'womens gradient leggings',
    'womens gray running shorts',
    'womens gray tights',
    'womens green compression pants',
    'womens guide pants',
    'womens guide pro pants',
    'womens gym apparel',
    'womens gym lycra',
    'womens gym pants with pockets',
    'womens gym tights',
    'womens head gloves',
    'womens high rise capris',
    'womens high stockings',
    'womens high waist capris',
    'womens high waist tights',
    'womens high waisted compression leggings',
    'womens hiking apparel',
    'womens hiking capris',
    'womens hiking tights',
    'womens hip pack',
    'womens insulated pants long',
    'womens insulated pants short',
    'womens jodphurs',
    'womens jogging capris',
    'womens knee length jacket',
    'womens knee length leggings',
    'womens knee length shorts',
    'womens knee length underwear',
    'womens knee length workout pants',
    'womens knee support',
    'womens legging shorts',
    'womens leggings 3 4 length',
    'womens leggings black workout',
    'womens leggings capri',
    'womens leggings full length',
    'womens leggings gym',
    'womens leggings knee length',
    'womens leggings low rise',
    'womens leggings lululemon',
    'womens leggings mid rise',
    'womens leggings running',
    'womens leggings sheer',
    'womens leggings short length',
    'womens leggings shorts',
    'womens leggings soft',
    'womens leggings tall long',
    'womens leggings tights',
    'womens leggings tops',
    'womens leggings under',
    'womens leggings white',
    'womens leggings white capri',
    'womens leggings with designs',
    'womens leggings with skirt',
    'womens leggings workout',
    'womens leggings workout high waisted',
    'womens lightweight capris',
    'womens lightweight trail pants',
    'womens lined running pants',
    'womens long compression shorts',
    'womens long running pants',
    'womens long running shorts',
    'womens long shorts to the knees',
    'womens long underwear leggings',
    'womens long underwear pants',
    'womens loose running pants',
    'womens lothes',
    'womens low rise capris',
    'womens low rise leggings',
    'womens lululemon leggings',
    'womens martial arts',
    'womens mid thigh bike shorts',
    'womens muscle building',
    'womens muscle tops',
    'womens nike air',
    'womens nike high tops',
    'womens nike pros',
    'womens nike therma fit pants',
    'womens nike therma pant',
    'womens nike therma training pants',
    'women''s nike tights',
    'womens nike training tights',
    'womens outdoor running gear',
    'womens panta clearance',
    'womens pants capri',
    'womens pants cold weather',
    'womens pants for the cold',
    'womens pants size 20',
    'womens pants under armour',
    'womens performance pants',
    'womens plus nike pants',
    'womens plus size capri leggings',
    'womens plus size cold gear',
    'womens plus size compression leggings',
    'womens plus size compression tights',
    'womens plus size stockings',
    'womens post workout',
    'womens puff vest',
    'womens puma tights',
    'womens quick dry capri pants'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
